$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.928.10"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "2.297.05"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'507.53"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'129.84"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "2.324.60"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("D10").Value = "'0.0980"
$ws.Range("E10").Value = "  +2.20%  "
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("E12").Value = "  +7.74%  "
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("D14").Value = "'23.88"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").Value = "2.707.46"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").Value = "54.897.66"
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("E17").Value = "  +1.38%  "
$ws.Range("D18").Value = "2.282.87"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Value = "'10.73"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").Value = "'4.19"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'6.68"
$ws.Range("E21").Value = "  +4.30%  "
$ws.Range("D22").Value = "'311.28"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'60.35"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").Value = "'0.151"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'7.53"
$ws.Range("E27").Value = "  +2.45%  "
$ws.Range("D28").Value = "'172.98"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").Value = "0.0₃0710"
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "'1.14"
$ws.Range("E32").Value = "  +4.28%  "
$ws.Range("D33").Value = "'18.10"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").Value = "'1.23"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").Value = "'0.917"
$ws.Range("E37").Value = "  -5.23%  "
$ws.Range("D38").Value = "'3.91"
$ws.Range("E38").Value = "  +3.27%  "
$ws.Range("D39").Value = "'36.74"
$ws.Range("E39").Value = "  +1.89%  "
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("E41").Value = "  +0.75%  "
$ws.Range("D42").Value = "'135.62"
$ws.Range("E42").Value = "  +7.69%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").Value = "'4.92"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").Value = "'259.79"
$ws.Range("E45").Value = "  +6.76%  "
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").Value = "'0.0913"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").Value = "'0.554"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'0.378"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("D50").Value = "'0.0210"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("E51").Value = "  +0.40%  "

Write-Host "Updated cryptos list"